$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the two newly-appearing "D" cells (existing shared string "modifier")
$ws.Range("D88").Value = "modifier"
$ws.Range("D97").Value = "modifier"

# Add a per-row (non-array) formula in column G for every data row (2-119):
# G{n} = IF(C{n}=E{n}, "-", "modifier")
for ($r = 2; $r -le 119; $r++) {
    $ws.Cells.Item($r, 7).Formula = "=IF(C$r=E$r, ""-"", ""modifier"")"
}

# Clear the autofilter criteria (colId 3 = "Action à faire") so every row is
# shown again, while leaving the autofilter buttons themselves in place.
$ws.ShowAllData()

# Unhide all previously hidden data rows (2-119)
$ws.Rows("2:119").Hidden = $false

# Update the sheet view: scroll position, zoom, active selection
$ws.Range("F7").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.Zoom = 95
